# Actualización desde MV -datos-
# Appends the next 5 daily observations (29-10-2021 .. 02-11-2021) to the
# bottom of the "Liquidez en pesos con incentivo al credito 2021 - Diaria"
# sheet, mirroring how the rest of the "Serie" column (A) is filled: plain
# text dates (dd-mm-yyyy) in column A, plain numbers in B:G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows, in order: Serie (date as text), then the 6 numeric columns.
$newRows = @(
    @("29-10-2021", 2104, 2625, 15180, 2650, 4588, 8815),
    @("30-10-2021", 2104, 2625, 15180, 2650, 4588, 8815),
    @("31-10-2021", 2104, 2625, 15180, 2650, 4588, 8815),
    @("01-11-2021", 2104, 2625, 15180, 2650, 4588, 8815),
    @("02-11-2021", 2090, 2607, 15079, 2632, 4558, 8756)
)

$startRow  = 303
# Scratch cell, far outside the used range, used only to force column A's
# values to be stored as literal text instead of being auto-parsed into
# date serials (Excel would otherwise silently convert strings such as
# "01-11-2021" / "02-11-2021" into dates because day <= 12 makes them look
# like valid mm-dd-yyyy dates). We format the scratch cell as Text, type
# the value there, copy its already-text value into the real cell with
# PasteSpecial (values only, so the destination keeps its original/default
# style), then clear the scratch cell again.
$scratchRow = 1048500
$scratch = $ws.Cells.Item($scratchRow, 1)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row    = $startRow + $i
    $values = $newRows[$i]

    # Column A: write through the text-formatted scratch cell so the date
    # string is preserved verbatim.
    $scratch.NumberFormat = "@"
    $scratch.Value = $values[0]
    $scratch.Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)
    $scratch.Clear()

    # Columns B:G: plain numeric values.
    for ($c = 1; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
